$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename Sheet3 -> "changeCurrentData"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws3.Name = "changeCurrentData"

# ---------------------------------------------------------------------------
# 2. changeCurrentData (old Sheet3): rebuild with new "I (A)" current column
#    and reordered B*L / d(B*L) / V / dV columns.
#    Final layout: A=I (A), B=B*L, C=d(B*L), D=V, E=dV
# ---------------------------------------------------------------------------
$ws3.Range("A1:E14").Clear()

$ws3.Cells.Item(1,1).Value = "I (A)"
$ws3.Cells.Item(1,2).Value = "B*L"
$ws3.Cells.Item(1,3).Value = "d(B*L)"
$ws3.Cells.Item(1,4).Value = "V"
$ws3.Cells.Item(1,5).Value = "dV"

$currents = @(-3,-2.5,-2,-1.5,-1,-0.5,0,0.5,1,1.5,2,2.5,3)
$bl       = @(-323,-269,-215,-161,-108,-54,0,54,108,161,215,269,323)
$dbl      = @(12,10,8,6,4,2,0,2,4,6,8,10,12)
$vvals    = @("2.8220000000000001","2.887","2.9510000000000001","3.016","3.0819999999999999","3.145","3.2109999999999999","3.2719999999999998","3.3380000000000001","3.4039999999999999","3.468","3.5339999999999998","3.5990000000000002")
$dvvals   = @("1.296E-10","5.4699999999999997E-10","3.7220000000000002E-10","3.7749999999999999E-11","8.98E-10","9.034E-10","1.184E-9","1.161E-10","2.0510000000000001E-9","2.133E-9","6.9580000000000002E-10","1.1160000000000001E-9","1.045E-9")

for ($i = 0; $i -lt 13; $i++) {
    $r = $i + 2
    $ws3.Cells.Item($r,1).Value = $currents[$i]
    $ws3.Cells.Item($r,2).Value = $bl[$i]
    $ws3.Cells.Item($r,3).Value = $dbl[$i]
    $ws3.Cells.Item($r,4).Value = [double]$vvals[$i]
    $ws3.Cells.Item($r,5).Value = [double]$dvvals[$i]
    $ws3.Cells.Item($r,5).NumberFormat = "0.00E+00"
}

$ws3.Range("A15").Select()

# ---------------------------------------------------------------------------
# 3. Sheet4: drop the old c-value summary rows (7-10) now that the "I (A)"
#    data lives on changeCurrentData.
# ---------------------------------------------------------------------------
$ws4.Rows.Item(7).Resize(4).Delete()
$ws4.Range("A3").Select()

# ---------------------------------------------------------------------------
# 4. Sheet2 cosmetic updates: column widths + becomes the active tab.
#    (ColumnWidth is stored/rounded to the nearest 1/6 character width by
#    the host, so 16.75/17.25 land on their nearest representable width.)
# ---------------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 16.75
$ws2.Columns.Item(2).ColumnWidth = 17.25
$ws2.Range("A6").Select()

# ---------------------------------------------------------------------------
# 5. Sheet1 selection unchanged; make Sheet2 the active sheet/tab.
# ---------------------------------------------------------------------------
$ws2.Activate()
